$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Censys Director Blog", "Counting Everyone Once, Only Once and in the Right Place", "https://www.census.gov/newsroom/blogs/director/2018/11/counting_everyoneon.html", "2021-10-08", "U.S. Census Bureau staff presented to the National Advisory Committee (NAC) and gave an update on how we plan to efficiently and effectively reach hard-to-count communities."),
    @("Censys Director Blog", "Thank You, Providence", "https://www.census.gov/newsroom/blogs/director/2018/09/thank_you_providenc.html", "2021-10-08", "Last week I had the opportunity to visit Providence, Rhode Island, where the 2018 Census Test is wrapping up."),
    @("Censys Director Blog", "Planning and Testing for the 2020 Census in Harris County, Texas", "https://www.census.gov/newsroom/blogs/director/2016/03/planning-and-testing-for-the-2020-census-in-harris-county-texas.html", "2021-10-08", "Today I am visiting Harris County, Texas, one of two sites now taking part in the 2016 Census Test that will help us prepare for the decennial census in 2020."),
    @("Censys Director Blog", "Measuring Income, Poverty, and Health Insurance Coverage", "https://www.census.gov/newsroom/blogs/director/2015/08/measuring-income-poverty-and-health-insurance-coverage.html", "2021-10-08", "Today, the U.S. Census Bureau hosted a webcast on income, poverty, and health insurance coverage statistics."),
    @("Censys Director Blog", "Transparency, Participation and Collaboration at the U.S. Census Bureau", "https://www.census.gov/newsroom/blogs/director/2015/05/transparency-participation-and-collaboration-at-the-u-s-census-bureau.html", "2022-03-22", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "2015 Census Test Starts Today in the Savannah, Ga. Area", "https://www.census.gov/newsroom/blogs/director/2015/03/2015-census-test-starts-today-in-the-savannah-ga-area-2.html", "2022-03-22", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Happy 25th Anniversary, TIGER", "https://www.census.gov/newsroom/blogs/director/2014/11/happy-25th-anniversary-tiger.html", "2022-03-22", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Introducing Census PoP Quiz, the Census Bureau�s Latest Mobile App", "https://www.census.gov/newsroom/blogs/director/2014/08/introducing-census-pop-quiz-the-census-bureaus-latest-mobile-app.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "224 Years of Measuring America's People, Places, and Economy", "https://www.census.gov/newsroom/blogs/director/2014/07/celebrating-224-years-of-measuring-americas-people-places-and-economy.html", "2021-12-16", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Collecting Reliable, Timely and Local Census Data", "https://www.census.gov/newsroom/blogs/director/2014/06/collecting-reliable-timely-and-local-census-data.html", "2022-03-22", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Know the Population, Every Minute, Every Day", "https://www.census.gov/newsroom/blogs/director/2013/04/know-the-population-every-minute-every-day.html", "2022-03-22", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "The Times, They Are a-Changin�", "https://www.census.gov/newsroom/blogs/director/2012/08/the-times-they-are-a-changin.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "And Now, for Something a Little Different . . .", "https://www.census.gov/newsroom/blogs/director/2012/06/and-now-for-something-a-little-different.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "How Good was the 2010 Census?, Part 3", "https://www.census.gov/newsroom/blogs/director/2012/05/how-good-was-the-2010-census-part-3.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Looking Forward, Looking Back", "https://www.census.gov/newsroom/blogs/director/2012/04/looking-forward-looking-back.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Genealogists, Start Your Engines!", "https://www.census.gov/newsroom/blogs/director/2012/03/genealogists-start-your-engines.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "National Statistical Offices: Independent, Alike, Real-Time Actions", "https://www.census.gov/newsroom/blogs/director/2012/02/national-statistical-offices-independent-identical-simultaneous-actions-thousands-of-miles-apart.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "The Future of Producing Social and Economic Statistical Information, Part III", "https://www.census.gov/newsroom/blogs/director/2011/10/the-future-of-producing-social-and-economic-statistical-information-part-iii.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "The Future of Producing Social and Economic Statistical Information, Part II", "https://www.census.gov/newsroom/blogs/director/2011/09/the-future-of-producing-social-and-economic-statistical-information-part-ii.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "The Future of Producing Social and Economic Statistical Information, Part I", "https://www.census.gov/newsroom/blogs/director/2011/09/the-future-of-producing-social-and-economic-statistical-information-part-i.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Quality in a Census, Some Overview Thoughts", "https://www.census.gov/newsroom/blogs/director/2010/09/quality-in-a-census-some-overview-thoughts.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Measuring Quality in a Census, Part 4", "https://www.census.gov/newsroom/blogs/director/2010/08/measuring-quality-in-a-census-part-4.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Measuring Quality in a Census, Part 3", "https://www.census.gov/newsroom/blogs/director/2010/07/measuring-quality-in-a-census-part-3.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Repairing a Problem, Connecting With Local Leaders", "https://www.census.gov/newsroom/blogs/director/2010/07/repairing-a-problem-connecting-with-local-leaders.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Measuring Quality in a Census, Part 2", "https://www.census.gov/newsroom/blogs/director/2010/07/measuring-quality-in-a-census-part-2.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Quality in a Census, Part 1", "https://www.census.gov/newsroom/blogs/director/2010/07/quality-in-a-census-part-1.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Progress Report, Five Weeks Into Nonresponse Followup", "https://www.census.gov/newsroom/blogs/director/2010/06/progress-report-five-weeks-into-nonresponse-followup.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Progress Report, Three Weeks Into Nonresponse Followup", "https://www.census.gov/newsroom/blogs/director/2010/05/progress-report-three-weeks-into-nonresponse-followup.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Progress Report, Two Weeks Into Nonresponse Followup", "https://www.census.gov/newsroom/blogs/director/2010/05/progress-report-two-weeks-into-nonresponse-followup.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Progress Report, One Week Into Nonresponse Followup", "https://www.census.gov/newsroom/blogs/director/2010/05/progress-report-one-week-into-nonresponse-followup.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "A Note to my 600,000 New Colleagues", "https://www.census.gov/newsroom/blogs/director/2010/05/a-note-to-my-600000-new-colleagues.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Academy Awards Show in LA", "https://www.census.gov/newsroom/blogs/director/2010/04/academy-awards-show-in-la.html", "2022-09-15", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "Our Partners, at Half-Time of The 2010 Census", "https://www.census.gov/newsroom/blogs/director/2010/04/our-partners-at-halftime-of-the-2010-census.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "So, How Are we Doing?", "https://www.census.gov/newsroom/blogs/director/2010/04/so-how-are-we-doing.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "10 Questions, 10 Minutes", "https://www.census.gov/newsroom/blogs/director/2010/04/10-questions-10-minutes.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "If You Got a Form, You�ll be Counted in The Right Location", "https://www.census.gov/newsroom/blogs/director/2010/03/if-you-got-a-form-youll-be-counted-in-the-right-location.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "So, How do You Handle Prisons?", "https://www.census.gov/newsroom/blogs/director/2010/03/so-how-do-you-handle-prisons.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy."),
    @("Censys Director Blog", "And so, it Began �", "https://www.census.gov/newsroom/blogs/director/2010/01/and-so-it-began.html", "2021-10-08", "Launched in 2009, this blog gives the director's point of view on everything from Census Bureau news to how we measure America's people, places and economy.")
)

$startRow = 628
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Style = "Normal"
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
